$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: 四方坪站
$ws.Cells.Item(8, 1).Value = 45995
$ws.Cells.Item(8, 2).Value = "四方坪站"
$ws.Cells.Item(8, 3).Value = 8580.26
$ws.Cells.Item(8, 4).Value = 7417.19
$ws.Cells.Item(8, 5).Value = 2850.22
$ws.Cells.Item(8, 6).Value = 392

# Row 9: 高岭站
$ws.Cells.Item(9, 1).Value = 45995
$ws.Cells.Item(9, 2).Value = "高岭站"
$ws.Cells.Item(9, 3).Value = 4940.75
$ws.Cells.Item(9, 4).Value = 4258.32
$ws.Cells.Item(9, 5).Value = 1413.67
$ws.Cells.Item(9, 6).Value = 170

# Copy styles from row 7 to rows 8-9 to keep formats consistent
$ws.Range("A7:F7").Copy()
$ws.Range("A8:F9").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("H10").Select()
